# Add a new diary row (row 7) to Sheet1 describing the new game idea
# ("add new game ifno by zhaoyj"), mirroring the layout/formatting that
# the existing row 6 uses.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- A7: "2012.4.14" ---------------------------------------------------
# Writing this text straight through Range.Value lets Excel's automatic
# date recognition kick in (it silently becomes a date serial number).
# To keep it as plain text we build it in a scratch cell via a formula
# (a quoted-string formula always yields a text result), copy it, and
# paste-special *values only* into A7. That keeps the literal string and
# does not create any new cell style.
$scratchA = $ws.Range("Z1000")
$scratchA.Formula = '="2012.4.14"'
$scratchA.Copy()
$ws.Range("A7").PasteSpecial(-4163)   # xlPasteValues
$scratchA.ClearContents()

# --- B7: game idea description -----------------------------------------
$ws.Range("B7").Value = "设计游戏《可爱糖果对对碰》思路，并上传"

# --- Match the formatting used by row 6 for these two new cells --------
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- D7: work duration, same as D6 --------------------------------------
$ws.Range("D7").Value = 2

# --- Restore the cursor/selection to the cell the author left active ---
$ws.Range("E6").Select() | Out-Null
